$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
